$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "eta^2" - copy formatting from the neighboring header
# cell (G1) so it gets the same bold/centered/bordered header style.
$ws.Range("H1").Value = "eta^2"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New effect size (eta^2) values in column H, rows 2-5. Force text storage
# (NumberFormat "@") so the numeric-looking values are written as strings,
# matching the convention used for every other stat value in this table.
# Then copy the formatting from the corresponding row's existing data cell
# (column G) so the new cells end up with the same (default/unstyled) look.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0.02"
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "0.25"
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "0.11"
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "0.11"
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
